# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 4150
$ws1.Range("F4").Value = 2388
$ws1.Range("F5").Value = 477
$ws1.Range("F7").Value = 39
$ws1.Range("F8").Value = 40
$ws1.Range("F11").Value = 109
$ws1.Range("F12").Value = 147
$ws1.Range("F13").Value = 1554
$ws1.Range("F15").Value = 3098
$ws1.Range("F16").Value = 212

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 4150
$ws4.Range("F4").Value = 2388
$ws4.Range("F5").Value = 477
$ws4.Range("F8").Value = 39
$ws4.Range("F9").Value = 40
$ws4.Range("F13").Value = 109
$ws4.Range("F14").Value = 147
$ws4.Range("F17").Value = 1554
$ws4.Range("F19").Value = 3098
$ws4.Range("F20").Value = 212
